$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "SS_Email" column (D) with header ---
$ws.Range("D1").Value2 = "SS_Email"
$ws.Range("D1").Font.Bold = $true

# --- Insert a brand-new row 2 (WCT / 55555555 / Closed) above the current data ---
$ws.Rows(2).Insert()
$ws.Range("A2").Value2 = 55555555
$ws.Range("C2").Value2 = "Closed"

# --- Fill in the SS_Email column for every data row with the hyperlinked address ---
$ws.Range("D2").Value2 = "tallen@mdsol.com"
$ws.Range("D3").Value2 = "tallen@mdsol.com"
$ws.Range("D4").Value2 = "tallen@mdsol.com"
$ws.Range("D5").Value2 = "tallen@mdsol.com"
$ws.Range("D6").Value2 = "tallen@mdsol.com"

$ws.Range("B2").Value2 = "WCT"

[void]$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:tallen@mdsol.com")
[void]$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:tallen@mdsol.com")
[void]$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:tallen@mdsol.com")
[void]$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:tallen@mdsol.com")
[void]$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:tallen@mdsol.com")

# Put the new row's font back to the regular (non-bold) style that was
# copied down from the header row during the insert.
$ws.Range("A2:C2").Font.Bold = $false
$ws.Range("A2:C2").Font.ThemeColor = 1

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 9.666666666666666
$ws.Columns("D").ColumnWidth = 15.498697916666666

# --- Selection / active cell ---
[void]$ws.Range("E2").Select()

Write-Output "done"
